$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.740.01"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "1.600.95"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'211.80"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "'19.59"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "1.825.86"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.602.73"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "'65.00"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'208.70"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'7.16"
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("D24").Value = "'143.61"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").Value = "1.280.96"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").Value = "'1.23"
$ws.Range("E35").Value = "  +16.65%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  -4.38%  "
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("D39").Value = "'0.824"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "'0.777"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").Value = "'62.63"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "1.738.19"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "'90.31"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "'0.103"
$ws.Range("E47").Value = "  +2.13%  "
$ws.Range("D48").Value = "'0.0512"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").Value = "'7.53"
$ws.Range("E49").Value = "  +3.46%  "
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  +1.78%  "
